# Add three new per-polymer "sink" result sheets (PA, PET, PS) after the
# existing Polymer_Results sheet, each holding a tiny Macro/Micro "Total Sum"
# summary table, mirroring the header style used on the first sheet.

$wb = $excel.ActiveWorkbook

# Template sheet whose A1:B1 header formatting (bold font + border + center
# alignment, style index "1" in styles.xml) we want the new sheets to reuse
# rather than re-creating a near-duplicate style resource.
$templateSheet = $wb.Worksheets.Item(1)

function Add-PolymerSinkSheet($SheetName, $MacroValue, $MicroValue) {

    # Insert the new sheet immediately after the current last tab so the
    # three end up in order right after Polymer_Results.
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $SheetName

    $ws.Range("A1").Value = " "
    $ws.Range("B1").Value = "Total Sum"
    $ws.Range("A2").Value = "Macro"
    $ws.Range("B2").Value = $MacroValue
    $ws.Range("A3").Value = "Micro"
    $ws.Range("B3").Value = $MicroValue

    # Copy the header cell formatting from Polymer_Results!A1:B1 onto row 1
    # of the new sheet so it reuses the existing bold/border/center style
    # instead of creating a brand-new font/xf pair.
    $templateSheet.Range("A1:B1").Copy()
    $ws.Range("A1:B1").PasteSpecial(-4122)

    # Match the 1"/1"/0.75"/0.75"... margins (in points: 72/72/54/54/36/36)
    # used by the rest of the workbook's sheets.
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

Add-PolymerSinkSheet "PA"  0.1302159545613743 86.02591692139549
Add-PolymerSinkSheet "PET" 432.9431841512236  266.268074969384
Add-PolymerSinkSheet "PS"  13.36698075189867  39.45174307950305
